$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the three "result" columns with the real astrology remedies
$ws.Range("D1").Value = "Good will happen"
$ws.Range("D2").Value = "Follow remedy 1"
$ws.Range("D3").Value = "Follow remedy 2"

# 2. Column A was using a distinct placeholder font; bring it (and the rest of
#    the used range) onto the same font already used by the data columns
#    (B/C/D) by copying that cell's format across the whole table.
$ws.Range("B1").Copy()
$ws.Range("A1:D6").PasteSpecial(-4122)

# 3. The sheet had several trailing, completely empty placeholder rows
#    (7-12) below the real 6-row table - drop them.
$ws.Range("A7:A12").EntireRow.Delete()

# 4. Center every cell in the (now 6-row) table.
$ws.Range("A1:D6").HorizontalAlignment = -4108
